# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Home row (row 2) with new offensive stats ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 552
$wsOff.Range("C2").Value = 389
$wsOff.Range("D2").Value = 136
$wsOff.Range("E2").Value = 56

# --- DEF sheet: update Home row (row 2) with new defensive stats ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 540
$wsDef.Range("C2").Value = 349
$wsDef.Range("D2").Value = 76
$wsDef.Range("E2").Value = 32
$wsDef.Range("F2").Value = 7
$wsDef.Range("G2").Value = 6
